$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "530×8=4240" "740×4=2960"
Replace-Text "249×4=996" "733×4=2932"
Replace-Text "477×6=2862" "853×7=5971"
Replace-Text "788×5=3940" "871×9=7839"
Replace-Text "988×5=4940" "321×4=1284"
Replace-Text "546×6=3276" "553×3=1659"
Replace-Text "120×3=360" "803×5=4015"
Replace-Text "551×2=1102" "965×6=5790"
Replace-Text "634×5=3170" "361×3=1083"
Replace-Text "996×5=4980" "126×8=1008"
Replace-Text "648×9=5832" "996×4=3984"
Replace-Text "631×4=2524" "418×4=1672"
Replace-Text "597×7=4179" "710×6=4260"
Replace-Text "974×7=6818" "166×9=1494"
Replace-Text "580×8=4640" "486×4=1944"
Replace-Text "175×2=350" "654×5=3270"
Replace-Text "250×6=1500" "538×6=3228"
Replace-Text "576×8=4608" "486×9=4374"
Replace-Text "141×3=423" "695×5=3475"
Replace-Text "390×5=1950" "147×6=882"
Replace-Text "863×8=6904" "611×3=1833"
Replace-Text "672×2=1344" "678×2=1356"
Replace-Text "813×5=4065" "435×6=2610"
Replace-Text "415×8=3320" "483×7=3381"
Replace-Text "378×7=2646" "186×3=558"
